# Update Plau-Plaur LR-pair stats following Dr Hou advice
# (recompute ligand/receptor-expressing-cell counts and all derived stats)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 22.556342
$ws.Range("H2").Value = 45.112684
$ws.Range("I2").Value = 0.0667629019027735
$ws.Range("J2").Value = 0.04665728030990886
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 14.4437745
$ws.Range("N2").Value = 28.887549
$ws.Range("O2").Value = 0.10117841194680531
$ws.Range("P2").Value = 0.07108478932534294
$ws.Range("Q2").Value = 325.798717392879
$ws.Range("R2").Value = 1303.194869571516
$ws.Range("S2").Value = 0.006754964391482969
$ws.Range("T2").Value = 0.0033166229413233426

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 22.556342
$ws.Range("H3").Value = 45.112684
$ws.Range("I3").Value = 0.0667629019027735
$ws.Range("J3").Value = 0.04665728030990886
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.722763333333333
$ws.Range("N3").Value = 20.16829
$ws.Range("O3").Value = 0.04709285083070758
$ws.Range("P3").Value = 0.04962894725691061
$ws.Range("Q3").Value = 151.64094893172665
$ws.Range("R3").Value = 909.8456935903599
$ws.Range("S3").Value = 0.003144055380332476
$ws.Range("T3").Value = 0.0023155517036513606

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 22.556342
$ws.Range("H4").Value = 45.112684
$ws.Range("I4").Value = 0.0667629019027735
$ws.Range("J4").Value = 0.04665728030990886
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.33817233333334
$ws.Range("N4").Value = 139.014517
$ws.Range("O4").Value = 0.3245981643651427
$ws.Range("P4").Value = 0.342078784673262
$ws.Range("Q4").Value = 1045.219662805605
$ws.Range("R4").Value = 6271.317976833629
$ws.Range("S4").Value = 0.021671115405330375
$ws.Range("T4").Value = 0.01596046574457334

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("G5").Value = 22.556342
$ws.Range("H5").Value = 45.112684
$ws.Range("I5").Value = 0.0667629019027735
$ws.Range("J5").Value = 0.04665728030990886
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 54.450333
$ws.Range("N5").Value = 163.350999
$ws.Range("O5").Value = 0.3814237215427815
$ws.Range("P5").Value = 0.4019645747723113
$ws.Range("Q5").Value = 1228.200333161886
$ws.Range("R5").Value = 7369.201998971316
$ws.Range("S5").Value = 0.025464954504751513
$ws.Range("T5").Value = 0.018754573839805045

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 22.556342
$ws.Range("H6").Value = 45.112684
$ws.Range("I6").Value = 0.0667629019027735
$ws.Range("J6").Value = 0.04665728030990886
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.359316999999999
$ws.Range("N6").Value = 40.077951
$ws.Range("O6").Value = 0.09358180430980552
$ws.Range("P6").Value = 0.09862147541234521
$ws.Range("Q6").Value = 301.337323138414
$ws.Range("R6").Value = 1808.023938830484
$ws.Range("S6").Value = 0.0062477928210200925
$ws.Range("T6").Value = 0.004601409822890575

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("G7").Value = 22.556342
$ws.Range("H7").Value = 45.112684
$ws.Range("I7").Value = 0.0667629019027735
$ws.Range("J7").Value = 0.04665728030990886
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 7.441137
$ws.Range("N7").Value = 14.882274
$ws.Range("O7").Value = 0.05212504700475731
$ws.Range("P7").Value = 0.03662142855982794
$ws.Range("Q7").Value = 167.84483104085402
$ws.Range("R7").Value = 671.3793241634161
$ws.Range("S7").Value = 0.00348001939985607
$ws.Range("T7").Value = 0.0017086562576651939

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 136.26068366666667
$ws.Range("H8").Value = 408.782051
$ws.Range("I8").Value = 0.40330824283664957
$ws.Range("J8").Value = 0.4227781867105593
$ws.Range("K8").Value = 2
$ws.Range("M8").Value = 14.4437745
$ws.Range("N8").Value = 28.887549
$ws.Range("O8").Value = 0.10117841194680531
$ws.Range("P8").Value = 0.07108478932534294
$ws.Range("Q8").Value = 1968.1185880971664
$ws.Range("R8").Value = 11808.711528583
$ws.Range("S8").Value = 0.040806087535268716
$ws.Range("T8").Value = 0.030053098333670612

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 136.26068366666667
$ws.Range("H9").Value = 408.782051
$ws.Range("I9").Value = 0.40330824283664957
$ws.Range("J9").Value = 0.4227781867105593
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.722763333333333
$ws.Range("N9").Value = 20.16829
$ws.Range("O9").Value = 0.04709285083070758
$ws.Range("P9").Value = 0.04962894725691061
$ws.Range("Q9").Value = 916.0483279291989
$ws.Range("R9").Value = 8244.43495136279
$ws.Range("S9").Value = 0.01899293491870113
$ws.Range("T9").Value = 0.020982036329630648

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 136.26068366666667
$ws.Range("H10").Value = 408.782051
$ws.Range("I10").Value = 0.40330824283664957
$ws.Range("J10").Value = 0.4227781867105593
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.33817233333334
$ws.Range("N10").Value = 139.014517
$ws.Range("O10").Value = 0.3245981643651427
$ws.Range("P10").Value = 0.342078784673262
$ws.Range("Q10").Value = 6314.0710420038195
$ws.Range("R10").Value = 56826.63937803438
$ws.Range("S10").Value = 0.13091311529810767
$ws.Range("T10").Value = 0.14462344829631357

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 136.26068366666667
$ws.Range("H11").Value = 408.782051
$ws.Range("I11").Value = 0.40330824283664957
$ws.Range("J11").Value = 0.4227781867105593
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 54.450333
$ws.Range("N11").Value = 163.350999
$ws.Range("O11").Value = 0.3814237215427815
$ws.Range("P11").Value = 0.4019645747723113
$ws.Range("Q11").Value = 7419.439600457661
$ws.Range("R11").Value = 66774.95640411896
$ws.Range("S11").Value = 0.1538313309116347
$ws.Range("T11").Value = 0.1699418540441188

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 136.26068366666667
$ws.Range("H12").Value = 408.782051
$ws.Range("I12").Value = 0.40330824283664957
$ws.Range("J12").Value = 0.4227781867105593
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 13.359316999999999
$ws.Range("N12").Value = 40.077951
$ws.Range("O12").Value = 0.09358180430980552
$ws.Range("P12").Value = 0.09862147541234521
$ws.Range("Q12").Value = 1820.3496677397222
$ws.Range("R12").Value = 16383.147009657501
$ws.Range("S12").Value = 0.03774231305767087
$ws.Range("T12").Value = 0.04169500854555131

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 136.26068366666667
$ws.Range("H13").Value = 408.782051
$ws.Range("I13").Value = 0.40330824283664957
$ws.Range("J13").Value = 0.4227781867105593
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 7.441137
$ws.Range("N13").Value = 14.882274
$ws.Range("O13").Value = 0.05212504700475731
$ws.Range("P13").Value = 0.03662142855982794
$ws.Range("Q13").Value = 1013.934414877329
$ws.Range("R13").Value = 6083.606489263974
$ws.Range("S13").Value = 0.021022461115266435
$ws.Range("T13").Value = 0.01548274116127434

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 60.036368
$ws.Range("H14").Value = 180.109104
$ws.Range("I14").Value = 0.17769734770747891
$ws.Range("J14").Value = 0.18627579222939888
$ws.Range("K14").Value = 2
$ws.Range("M14").Value = 14.4437745
$ws.Range("N14").Value = 28.887549
$ws.Range("O14").Value = 0.10117841194680531
$ws.Range("P14").Value = 0.07108478932534294
$ws.Range("Q14").Value = 867.151761191016
$ws.Range("R14").Value = 5202.910567146096
$ws.Range("S14").Value = 0.017979135448202
$ws.Range("T14").Value = 0.013241375447038172

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 60.036368
$ws.Range("H15").Value = 180.109104
$ws.Range("I15").Value = 0.17769734770747891
$ws.Range("J15").Value = 0.18627579222939888
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 6.722763333333333
$ws.Range("N15").Value = 20.16829
$ws.Range("O15").Value = 0.04709285083070758
$ws.Range("P15").Value = 0.04962894725691061
$ws.Range("Q15").Value = 403.6102934569067
$ws.Range("R15").Value = 3632.49264111216
$ws.Range("S15").Value = 0.008368274688600682
$ws.Range("T15").Value = 0.009244671467792075

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 60.036368
$ws.Range("H16").Value = 180.109104
$ws.Range("I16").Value = 0.17769734770747891
$ws.Range("J16").Value = 0.18627579222939888
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 46.33817233333334
$ws.Range("N16").Value = 139.014517
$ws.Range("O16").Value = 0.3245981643651427
$ws.Range("P16").Value = 0.342078784673262
$ws.Range("Q16").Value = 2781.9755666514193
$ws.Range("R16").Value = 25037.78009986277
$ws.Range("S16").Value = 0.05768023287840216
$ws.Range("T16").Value = 0.06372099661988184

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 60.036368
$ws.Range("H17").Value = 180.109104
$ws.Range("I17").Value = 0.17769734770747891
$ws.Range("J17").Value = 0.18627579222939888
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 54.450333
$ws.Range("N17").Value = 163.350999
$ws.Range("O17").Value = 0.3814237215427815
$ws.Range("P17").Value = 0.4019645747723113
$ws.Range("Q17").Value = 3269.0002297105443
$ws.Range("R17").Value = 29421.002067394897
$ws.Range("S17").Value = 0.06777798367086825
$ws.Range("T17").Value = 0.07487626961386572

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 60.036368
$ws.Range("H18").Value = 180.109104
$ws.Range("I18").Value = 0.17769734770747891
$ws.Range("J18").Value = 0.18627579222939888
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 13.359316999999999
$ws.Range("N18").Value = 40.077951
$ws.Range("O18").Value = 0.09358180430980552
$ws.Range("P18").Value = 0.09862147541234521
$ws.Range("Q18").Value = 802.044871640656
$ws.Range("R18").Value = 7218.403844765904
$ws.Range("S18").Value = 0.01662923841953276
$ws.Range("T18").Value = 0.018370793463266786

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 60.036368
$ws.Range("H19").Value = 180.109104
$ws.Range("I19").Value = 0.17769734770747891
$ws.Range("J19").Value = 0.18627579222939888
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 7.441137
$ws.Range("N19").Value = 14.882274
$ws.Range("O19").Value = 0.05212504700475731
$ws.Range("P19").Value = 0.03662142855982794
$ws.Range("Q19").Value = 446.73883927041607
$ws.Range("R19").Value = 2680.433035622496
$ws.Range("S19").Value = 0.009262482601873042
$ws.Range("T19").Value = 0.006821685617554281

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 67.73010366666666
$ws.Range("H20").Value = 203.190311
$ws.Range("I20").Value = 0.20046948512140605
$ws.Range("J20").Value = 0.2101472680407257
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 14.4437745
$ws.Range("N20").Value = 28.887549
$ws.Range("O20").Value = 0.10117841194680531
$ws.Range("P20").Value = 0.07108478932534294
$ws.Range("Q20").Value = 978.2783442229564
$ws.Range("R20").Value = 5869.670065337739
$ws.Range("S20").Value = 0.020283184148377577
$ws.Range("T20").Value = 0.01493827427597136

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 67.73010366666666
$ws.Range("H21").Value = 203.190311
$ws.Range("I21").Value = 0.20046948512140605
$ws.Range("J21").Value = 0.2101472680407257
$ws.Range("K21").Value = 3
$ws.Range("M21").Value = 6.722763333333333
$ws.Range("N21").Value = 20.16829
$ws.Range("O21").Value = 0.04709285083070758
$ws.Range("P21").Value = 0.04962894725691061
$ws.Range("Q21").Value = 455.3334574931322
$ws.Range("R21").Value = 4098.00111743819
$ws.Range("S21").Value = 0.009440679558931128
$ws.Range("T21").Value = 0.01042938768177703

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("G22").Value = 67.73010366666666
$ws.Range("H22").Value = 203.190311
$ws.Range("I22").Value = 0.20046948512140605
$ws.Range("J22").Value = 0.2101472680407257
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 46.33817233333334
$ws.Range("N22").Value = 139.014517
$ws.Range("O22").Value = 0.3245981643651427
$ws.Range("P22").Value = 0.342078784673262
$ws.Range("Q22").Value = 3138.489215860532
$ws.Range("R22").Value = 28246.40294274479
$ws.Range("S22").Value = 0.06507202688163369
$ws.Range("T22").Value = 0.07188692205377768

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("G23").Value = 67.73010366666666
$ws.Range("H23").Value = 203.190311
$ws.Range("I23").Value = 0.20046948512140605
$ws.Range("J23").Value = 0.2101472680407257
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 54.450333
$ws.Range("N23").Value = 163.350999
$ws.Range("O23").Value = 0.3814237215427815
$ws.Range("P23").Value = 0.4019645747723113
$ws.Range("Q23").Value = 3687.926698774521
$ws.Range("R23").Value = 33191.34028897069
$ws.Range("S23").Value = 0.07646381707077195
$ws.Range("T23").Value = 0.08447175723755324

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("G24").Value = 67.73010366666666
$ws.Range("H24").Value = 203.190311
$ws.Range("I24").Value = 0.20046948512140605
$ws.Range("J24").Value = 0.2101472680407257
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 13.359316999999999
$ws.Range("N24").Value = 40.077951
$ws.Range("O24").Value = 0.09358180430980552
$ws.Range("P24").Value = 0.09862147541234521
$ws.Range("Q24").Value = 904.8279253258622
$ws.Range("R24").Value = 8143.451327932761
$ws.Range("S24").Value = 0.01876029612671889
$ws.Range("T24").Value = 0.02072503362804995

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("G25").Value = 67.73010366666666
$ws.Range("H25").Value = 203.190311
$ws.Range("I25").Value = 0.20046948512140605
$ws.Range("J25").Value = 0.2101472680407257
$ws.Range("K25").Value = 2
$ws.Range("M25").Value = 7.441137
$ws.Range("N25").Value = 14.882274
$ws.Range("O25").Value = 0.05212504700475731
$ws.Range("P25").Value = 0.03662142855982794
$ws.Range("Q25").Value = 503.988980407869
$ws.Range("R25").Value = 3023.9338824472143
$ws.Range("S25").Value = 0.010449481334972786
$ws.Range("T25").Value = 0.0076958931635964485

# Row 26
$ws.Range("E26").Value = 3
$ws.Range("G26").Value = 27.15284733333333
$ws.Range("H26").Value = 81.458542
$ws.Range("I26").Value = 0.08036776898028582
$ws.Range("J26").Value = 0.08424757054425056
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 14.4437745
$ws.Range("N26").Value = 28.887549
$ws.Range("O26").Value = 0.10117841194680531
$ws.Range("P26").Value = 0.07108478932534294
$ws.Range("Q26").Value = 392.18960391559295
$ws.Range("R26").Value = 2353.137623493558
$ws.Range("S26").Value = 0.008131483237133039
$ws.Range("T26").Value = 0.005988720803310019

# Row 27
$ws.Range("E27").Value = 3
$ws.Range("G27").Value = 27.15284733333333
$ws.Range("H27").Value = 81.458542
$ws.Range("I27").Value = 0.08036776898028582
$ws.Range("J27").Value = 0.08424757054425056
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 6.722763333333333
$ws.Range("N27").Value = 20.16829
$ws.Range("O27").Value = 0.04709285083070758
$ws.Range("P27").Value = 0.04962894725691061
$ws.Range("Q27").Value = 182.54216644813107
$ws.Range("R27").Value = 1642.8794980331797
$ws.Range("S27").Value = 0.0037847473561853674
$ws.Range("T27").Value = 0.004181118235063466

# Row 28
$ws.Range("E28").Value = 3
$ws.Range("G28").Value = 27.15284733333333
$ws.Range("H28").Value = 81.458542
$ws.Range("I28").Value = 0.08036776898028582
$ws.Range("J28").Value = 0.08424757054425056
$ws.Range("K28").Value = 3
$ws.Range("M28").Value = 46.33817233333334
$ws.Range("N28").Value = 139.014517
$ws.Range("O28").Value = 0.3245981643651427
$ws.Range("P28").Value = 0.342078784673262
$ws.Range("Q28").Value = 1258.2133190726904
$ws.Range("R28").Value = 11323.919871654214
$ws.Range("S28").Value = 0.026087230285122633
$ws.Range("T28").Value = 0.02881930654345214

# Row 29
$ws.Range("E29").Value = 3
$ws.Range("G29").Value = 27.15284733333333
$ws.Range("H29").Value = 81.458542
$ws.Range("I29").Value = 0.08036776898028582
$ws.Range("J29").Value = 0.08424757054425056
$ws.Range("K29").Value = 3
$ws.Range("M29").Value = 54.450333
$ws.Range("N29").Value = 163.350999
$ws.Range("O29").Value = 0.3814237215427815
$ws.Range("P29").Value = 0.4019645747723113
$ws.Range("Q29").Value = 1478.4815791981619
$ws.Range("R29").Value = 13306.334212783457
$ws.Range("S29").Value = 0.030654173536551126
$ws.Range("T29").Value = 0.03386453886941997

# Row 30
$ws.Range("E30").Value = 3
$ws.Range("G30").Value = 27.15284733333333
$ws.Range("H30").Value = 81.458542
$ws.Range("I30").Value = 0.08036776898028582
$ws.Range("J30").Value = 0.08424757054425056
$ws.Range("K30").Value = 3
$ws.Range("M30").Value = 13.359316999999999
$ws.Range("N30").Value = 40.077951
$ws.Range("O30").Value = 0.09358180430980552
$ws.Range("P30").Value = 0.09862147541234521
$ws.Range("Q30").Value = 362.7434949786046
$ws.Range("R30").Value = 3264.6914548074415
$ws.Range("S30").Value = 0.007520960829528766
$ws.Range("T30").Value = 0.008308619706979624

# Row 31
$ws.Range("E31").Value = 3
$ws.Range("G31").Value = 27.15284733333333
$ws.Range("H31").Value = 81.458542
$ws.Range("I31").Value = 0.08036776898028582
$ws.Range("J31").Value = 0.08424757054425056
$ws.Range("K31").Value = 2
$ws.Range("M31").Value = 7.441137
$ws.Range("N31").Value = 14.882274
$ws.Range("O31").Value = 0.05212504700475731
$ws.Range("P31").Value = 0.03662142855982794
$ws.Range("Q31").Value = 202.048056947418
$ws.Range("R31").Value = 1212.288341684508
$ws.Range("S31").Value = 0.004189173735764875
$ws.Range("T31").Value = 0.003085266386025336

# Row 32
$ws.Range("E32").Value = 2
$ws.Range("G32").Value = 24.1210785
$ws.Range("H32").Value = 48.242157
$ws.Range("I32").Value = 0.07139425345140621
$ws.Range("J32").Value = 0.04989390216515674
$ws.Range("K32").Value = 2
$ws.Range("M32").Value = 14.4437745
$ws.Range("N32").Value = 28.887549
$ws.Range("O32").Value = 0.10117841194680531
$ws.Range("P32").Value = 0.07108478932534294
$ws.Range("Q32").Value = 348.39941855079826
$ws.Range("R32").Value = 1393.597674203193
$ws.Range("S32").Value = 0.0072235571863410045
$ws.Range("T32").Value = 0.003546697524029439

# Row 33
$ws.Range("E33").Value = 2
$ws.Range("G33").Value = 24.1210785
$ws.Range("H33").Value = 48.242157
$ws.Range("I33").Value = 0.07139425345140621
$ws.Range("J33").Value = 0.04989390216515674
$ws.Range("K33").Value = 3
$ws.Range("M33").Value = 6.722763333333333
$ws.Range("N33").Value = 20.16829
$ws.Range("O33").Value = 0.04709285083070758
$ws.Range("P33").Value = 0.04962894725691061
$ws.Range("Q33").Value = 162.160302100255
$ws.Range("R33").Value = 972.9618126015299
$ws.Range("S33").Value = 0.0033621589279568026
$ws.Range("T33").Value = 0.0024761818389960216

# Row 34
$ws.Range("E34").Value = 2
$ws.Range("G34").Value = 24.1210785
$ws.Range("H34").Value = 48.242157
$ws.Range("I34").Value = 0.07139425345140621
$ws.Range("J34").Value = 0.04989390216515674
$ws.Range("K34").Value = 3
$ws.Range("M34").Value = 46.33817233333334
$ws.Range("N34").Value = 139.014517
$ws.Range("O34").Value = 0.3245981643651427
$ws.Range("P34").Value = 0.342078784673262
$ws.Range("Q34").Value = 1117.7266923988616
$ws.Range("R34").Value = 6706.360154393169
$ws.Range("S34").Value = 0.02317444361654621
$ws.Range("T34").Value = 0.01706764541526345

# Row 35
$ws.Range("E35").Value = 2
$ws.Range("G35").Value = 24.1210785
$ws.Range("H35").Value = 48.242157
$ws.Range("I35").Value = 0.07139425345140621
$ws.Range("J35").Value = 0.04989390216515674
$ws.Range("K35").Value = 3
$ws.Range("M35").Value = 54.450333
$ws.Range("N35").Value = 163.350999
$ws.Range("O35").Value = 0.3814237215427815
$ws.Range("P35").Value = 0.4019645747723113
$ws.Range("Q35").Value = 1313.4007566441405
$ws.Range("R35").Value = 7880.404539864843
$ws.Range("S35").Value = 0.027231461848203926
$ws.Range("T35").Value = 0.02005558116754853

# Row 36
$ws.Range("E36").Value = 2
$ws.Range("G36").Value = 24.1210785
$ws.Range("H36").Value = 48.242157
$ws.Range("I36").Value = 0.07139425345140621
$ws.Range("J36").Value = 0.04989390216515674
$ws.Range("K36").Value = 3
$ws.Range("M36").Value = 13.359316999999999
$ws.Range("N36").Value = 40.077951
$ws.Range("O36").Value = 0.09358180430980552
$ws.Range("P36").Value = 0.09862147541234521
$ws.Range("Q36").Value = 322.2411340633845
$ws.Range("R36").Value = 1933.446804380307
$ws.Range("S36").Value = 0.006681203055334154
$ws.Range("T36").Value = 0.0049206102456069625

# Row 37
$ws.Range("E37").Value = 2
$ws.Range("G37").Value = 24.1210785
$ws.Range("H37").Value = 48.242157
$ws.Range("I37").Value = 0.07139425345140621
$ws.Range("J37").Value = 0.04989390216515674
$ws.Range("K37").Value = 2
$ws.Range("M37").Value = 7.441137
$ws.Range("N37").Value = 14.882274
$ws.Range("O37").Value = 0.05212504700475731
$ws.Range("P37").Value = 0.03662142855982794
$ws.Range("Q37").Value = 179.4882497062545
$ws.Range("R37").Value = 717.952998825018
$ws.Range("S37").Value = 0.0037214288170241056
$ws.Range("T37").Value = 0.0018271859737123317

